# Correcting Relevance Markers Walker (2018) - Wolters (2018)
# Updates the td_sim_1 (column C) simulation results and the corresponding
# record_atd (column D) values, plus the recomputed average in C222.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 94
$ws.Range("D3").Value = 94
$ws.Range("C5").Value = 48
$ws.Range("D5").Value = 48
$ws.Range("C7").Value = 268
$ws.Range("D7").Value = 268
$ws.Range("C9").Value = 127
$ws.Range("D9").Value = 127
$ws.Range("C11").Value = 74
$ws.Range("D11").Value = 74
$ws.Range("C13").Value = 38
$ws.Range("D13").Value = 38
$ws.Range("C15").Value = 36
$ws.Range("D15").Value = 36
$ws.Range("C19").Value = 93
$ws.Range("D19").Value = 93
$ws.Range("C21").Value = 67
$ws.Range("D21").Value = 67
$ws.Range("C23").Value = 117
$ws.Range("D23").Value = 117
$ws.Range("C25").Value = 73
$ws.Range("D25").Value = 73
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("C29").Value = 54
$ws.Range("D29").Value = 54
$ws.Range("C31").Value = 151
$ws.Range("D31").Value = 151
$ws.Range("C33").Value = 129
$ws.Range("D33").Value = 129
$ws.Range("C35").Value = 58
$ws.Range("D35").Value = 58
$ws.Range("C37").Value = 10
$ws.Range("D37").Value = 10
$ws.Range("C39").Value = 116
$ws.Range("D39").Value = 116
$ws.Range("C41").Value = 414
$ws.Range("D41").Value = 414
$ws.Range("C43").Value = 900
$ws.Range("D43").Value = 900
$ws.Range("C45").Value = 516
$ws.Range("D45").Value = 516
$ws.Range("C47").Value = 80
$ws.Range("D47").Value = 80
$ws.Range("C49").Value = 62
$ws.Range("D49").Value = 62
$ws.Range("C51").Value = 1
$ws.Range("D51").Value = 1
$ws.Range("C53").Value = 91
$ws.Range("D53").Value = 91
$ws.Range("C55").Value = 76
$ws.Range("D55").Value = 76
$ws.Range("C57").Value = 81
$ws.Range("D57").Value = 81
$ws.Range("C59").Value = 134
$ws.Range("D59").Value = 134
$ws.Range("C62").Value = 209
$ws.Range("D62").Value = 209
$ws.Range("C63").Value = 5
$ws.Range("D63").Value = 5
$ws.Range("C65").Value = 1502
$ws.Range("D65").Value = 1502
$ws.Range("C67").Value = 137
$ws.Range("D67").Value = 137
$ws.Range("C69").Value = 23
$ws.Range("D69").Value = 23
$ws.Range("C71").Value = 124
$ws.Range("D71").Value = 124
$ws.Range("C73").Value = 35
$ws.Range("D73").Value = 35
$ws.Range("C75").Value = 16
$ws.Range("D75").Value = 16
$ws.Range("C77").Value = 332
$ws.Range("D77").Value = 332
$ws.Range("C79").Value = 318
$ws.Range("D79").Value = 318
$ws.Range("C81").Value = 14
$ws.Range("D81").Value = 14
$ws.Range("C83").Value = 114
$ws.Range("D83").Value = 114
$ws.Range("C85").Value = 211
$ws.Range("D85").Value = 211
$ws.Range("C87").Value = 451
$ws.Range("D87").Value = 451
$ws.Range("C89").Value = 978
$ws.Range("D89").Value = 978
$ws.Range("C91").Value = 766
$ws.Range("D91").Value = 766
$ws.Range("C93").Value = 21
$ws.Range("D93").Value = 21
$ws.Range("C95").Value = 644
$ws.Range("D95").Value = 644
$ws.Range("C97").Value = 69
$ws.Range("D97").Value = 69
$ws.Range("C99").Value = 32
$ws.Range("D99").Value = 32
$ws.Range("C101").Value = 28
$ws.Range("D101").Value = 28
$ws.Range("C103").Value = 19
$ws.Range("D103").Value = 19
$ws.Range("C105").Value = 153
$ws.Range("D105").Value = 153
$ws.Range("C107").Value = 47
$ws.Range("D107").Value = 47
$ws.Range("C109").Value = 350
$ws.Range("D109").Value = 350
$ws.Range("C111").Value = 1670
$ws.Range("D111").Value = 1670
$ws.Range("C113").Value = 152
$ws.Range("D113").Value = 152
$ws.Range("C115").Value = 71
$ws.Range("D115").Value = 71
$ws.Range("C117").Value = 11
$ws.Range("D117").Value = 11
$ws.Range("C119").Value = 1887
$ws.Range("D119").Value = 1887
$ws.Range("C121").Value = 57
$ws.Range("D121").Value = 57
$ws.Range("C123").Value = 165
$ws.Range("D123").Value = 165
$ws.Range("C125").Value = 1374
$ws.Range("D125").Value = 1374
$ws.Range("C126").Value = 125
$ws.Range("D126").Value = 125
$ws.Range("C128").Value = 321
$ws.Range("D128").Value = 321
$ws.Range("C130").Value = 8
$ws.Range("D130").Value = 8
$ws.Range("C132").Value = 1161
$ws.Range("D132").Value = 1161
$ws.Range("C134").Value = 25
$ws.Range("D134").Value = 25
$ws.Range("C136").Value = 24
$ws.Range("D136").Value = 24
$ws.Range("C138").Value = 64
$ws.Range("D138").Value = 64
$ws.Range("C140").Value = 15
$ws.Range("D140").Value = 15
$ws.Range("C142").Value = 98
$ws.Range("D142").Value = 98
$ws.Range("C144").Value = 68
$ws.Range("D144").Value = 68
$ws.Range("C146").Value = 349
$ws.Range("D146").Value = 349
$ws.Range("C148").Value = 30
$ws.Range("D148").Value = 30
$ws.Range("C150").Value = 50
$ws.Range("D150").Value = 50
$ws.Range("C152").Value = 193
$ws.Range("D152").Value = 193
$ws.Range("C154").Value = 79
$ws.Range("D154").Value = 79
$ws.Range("C156").Value = 77
$ws.Range("D156").Value = 77
$ws.Range("C158").Value = 59
$ws.Range("D158").Value = 59
$ws.Range("C160").Value = 263
$ws.Range("D160").Value = 263
$ws.Range("C162").Value = 41
$ws.Range("D162").Value = 41
$ws.Range("C164").Value = 1168
$ws.Range("D164").Value = 1168
$ws.Range("C166").Value = 424
$ws.Range("D166").Value = 424
$ws.Range("C168").Value = 1181
$ws.Range("D168").Value = 1181
$ws.Range("C170").Value = 555
$ws.Range("D170").Value = 555
$ws.Range("C172").Value = 2180
$ws.Range("D172").Value = 2180
$ws.Range("C174").Value = 738
$ws.Range("D174").Value = 738
$ws.Range("C176").Value = 185
$ws.Range("D176").Value = 185
$ws.Range("C178").Value = 46
$ws.Range("D178").Value = 46
$ws.Range("C180").Value = 18
$ws.Range("D180").Value = 350.5
$ws.Range("C181").Value = 613
$ws.Range("D181").Value = 613
$ws.Range("C183").Value = 4
$ws.Range("D183").Value = 4
$ws.Range("C185").Value = 145
$ws.Range("D185").Value = 145
$ws.Range("C187").Value = 356
$ws.Range("D187").Value = 356
$ws.Range("C189").Value = 29
$ws.Range("D189").Value = 29
$ws.Range("C191").Value = 7
$ws.Range("D191").Value = 7
$ws.Range("C193").Value = 293
$ws.Range("D193").Value = 293
$ws.Range("C195").Value = 6
$ws.Range("D195").Value = 6
$ws.Range("C197").Value = 348
$ws.Range("D197").Value = 348
$ws.Range("C199").Value = 65
$ws.Range("D199").Value = 65
$ws.Range("C201").Value = 49
$ws.Range("D201").Value = 49
$ws.Range("C203").Value = 154
$ws.Range("D203").Value = 154
$ws.Range("C205").Value = 52
$ws.Range("D205").Value = 52
$ws.Range("C207").Value = 70
$ws.Range("D207").Value = 70
$ws.Range("C209").Value = 33
$ws.Range("D209").Value = 33
$ws.Range("C211").Value = 483
$ws.Range("D211").Value = 483
$ws.Range("C213").Value = 786
$ws.Range("D213").Value = 786
$ws.Range("C215").Value = 979
$ws.Range("D215").Value = 979
$ws.Range("C217").Value = 527
$ws.Range("D217").Value = 527
$ws.Range("C219").Value = 55
$ws.Range("D219").Value = 55
$ws.Range("C221").Value = 43
$ws.Range("D221").Value = 43
$ws.Range("C222").Value = 276.0990990990991
